$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before the existing "Outstanding" column (N),
# pushing Outstanding -> O and Original -> Q (Variable Instalments column).
$ws3.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab and set its selection.
$ws3.Activate()
$ws3.Range("T14").Select()
